# Weekly driver report update for 2025-04-21
# Updates the "Bad Drivers" and "Good Drivers (Roaming > 99.8%)" tables on the
# active sheet of LSA_driver_summary.xlsx with this week's roaming numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) "Bad Drivers" table
# ---------------------------------------------------------------------------

# Refresh the Critical Minutes / Good Roaming % for the three existing rows
$ws.Range("C3").Value = 1527
$ws.Range("D3").Value = 51.6

$ws.Range("C4").Value = 3213
$ws.Range("D4").Value = 93.4

$ws.Range("C5").Value = 339
$ws.Range("D5").Value = 96.7

# A new bad driver showed up this week - make room for it just above "Totals:"
$ws.Rows.Item(6).Insert()

$ws.Range("A6").Value = "Intel(R) Wi-Fi 6E AX210 160MHz - 23.90.0.2"
$ws.Range("B6").Value = 4
$ws.Range("C6").Value = 39
$ws.Range("D6").Value = 98.9

# "Totals:" row is now row 7 - update the rolled-up counts
$ws.Range("B7").Value = 29
$ws.Range("C7").Value = 5118

# ---------------------------------------------------------------------------
# 2) "Good Drivers (Roaming > 99.8%)" table
# ---------------------------------------------------------------------------
# The header block is now at rows 13-14 (shifted down by the insert above).
# The table grew from 9 to 17 driver rows, so insert 8 more rows after the
# header, then copy the existing row formatting down before writing values.

for ($i = 0; $i -lt 8; $i++) {
    $ws.Rows.Item(15).Insert()
}

$ws.Range("A23:E23").Copy()
$ws.Range("A15:E22").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 15: Intel(R) Wi-Fi 6E AX210 160MHz - 22.0.1.5 (no vintage recorded)
$ws.Range("A15").Value = "Intel(R) Wi-Fi 6E AX210 160MHz - 22.0.1.5"
$ws.Range("B15").Value = 156943
$ws.Range("D15").Value = 100
$ws.Range("E15").ClearContents()

# Row 16: Intel(R) Wi-Fi 6E AX210 160MHz - 23.120.0.3
$ws.Range("A16").Value = "Intel(R) Wi-Fi 6E AX210 160MHz - 23.120.0.3"
$ws.Range("B16").Value = 34181
$ws.Range("D16").Value = 99.9
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "2025-02-05"

# Row 17: Intel(R) Wi-Fi 6 AX201 160MHz - 23.100.0.4
$ws.Range("A17").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 23.100.0.4"
$ws.Range("B17").Value = 449371
$ws.Range("D17").Value = 99.9
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "2024-11-10"

# Row 18: Intel(R) Wi-Fi 6E AX210 160MHz - 23.20.1.1
$ws.Range("A18").Value = "Intel(R) Wi-Fi 6E AX210 160MHz - 23.20.1.1"
$ws.Range("B18").Value = 14968
$ws.Range("D18").Value = 100
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "2023-12-19"

# Row 19: Intel(R) Wi-Fi 6E AX210 160MHz - 22.170.2.1
$ws.Range("A19").Value = "Intel(R) Wi-Fi 6E AX210 160MHz - 22.170.2.1"
$ws.Range("B19").Value = 19083
$ws.Range("D19").Value = 100
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "2022-08-30"

# Row 20: Intel(R) Wi-Fi 6E AX211 160MHz - 22.150.3.1
$ws.Range("A20").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 22.150.3.1"
$ws.Range("B20").Value = 10661
$ws.Range("D20").Value = 100
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "2022-08-29"

# Row 21: Intel(R) Wi-Fi 6E AX211 160MHz - 22.150.0.3
$ws.Range("A21").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 22.150.0.3"
$ws.Range("B21").Value = 14239
$ws.Range("D21").Value = 100
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "2022-05-23"

# Row 22: Intel(R) Wi-Fi 6E AX211 160MHz - 22.100.1.1
$ws.Range("A22").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 22.100.1.1"
$ws.Range("B22").Value = 265400
$ws.Range("D22").Value = 99.9
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "2022-05-01"

# Row 23: Intel(R) Wi-Fi 6E AX210 160MHz - 22.100.0.3
$ws.Range("A23").Value = "Intel(R) Wi-Fi 6E AX210 160MHz - 22.100.0.3"
$ws.Range("B23").Value = 12988
$ws.Range("D23").Value = 100
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "2022-05-01"

# Row 24: Intel(R) Wi-Fi 6E AX210 160MHz - 22.130.0.5
$ws.Range("A24").Value = "Intel(R) Wi-Fi 6E AX210 160MHz - 22.130.0.5"
$ws.Range("B24").Value = 18738
$ws.Range("D24").Value = 99.9
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "2022-03-14"

# Row 25: Intel(R) Wi-Fi 6E AX210 160MHz - 22.110.1.1
$ws.Range("A25").Value = "Intel(R) Wi-Fi 6E AX210 160MHz - 22.110.1.1"
$ws.Range("B25").Value = 42024
$ws.Range("D25").Value = 100
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "2022-01-01"

# Row 26: Intel(R) Wi-Fi 6 AX201 160MHz - 22.80.0.9
$ws.Range("A26").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.80.0.9"
$ws.Range("B26").Value = 77999
$ws.Range("D26").Value = 99.9
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "2021-08-18"

# Row 27: Intel(R) Wi-Fi 6E AX210 160MHz - 22.70.0.6
$ws.Range("A27").Value = "Intel(R) Wi-Fi 6E AX210 160MHz - 22.70.0.6"
$ws.Range("B27").Value = 15504
$ws.Range("D27").Value = 100
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "2021-06-28"

# Row 28: Intel(R) Wi-Fi 6 AX201 160MHz - 22.50.1.1
$ws.Range("A28").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.50.1.1"
$ws.Range("B28").Value = 34244
$ws.Range("D28").Value = 100
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "2021-04-27"

# Row 29: Intel(R) Wi-Fi 6 AX201 160MHz - 21.110.3.2
$ws.Range("A29").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.110.3.2"
$ws.Range("B29").Value = 59673
$ws.Range("D29").Value = 100
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "2020-08-05"

# Row 30: Intel(R) Wi-Fi 6 AX201 160MHz - 21.70.0.6
$ws.Range("A30").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.70.0.6"
$ws.Range("B30").Value = 113652
$ws.Range("D30").Value = 100
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "2020-01-06"

# Row 31: Intel(R) Wi-Fi 6 AX201 160MHz - 21.60.2.1
$ws.Range("A31").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.60.2.1"
$ws.Range("B31").Value = 56018
$ws.Range("D31").Value = 100
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "2019-12-14"
